# Mudou logo para laranja no REMOTO
# Set the "Elipse 3" (logo) shape's fill to the theme's accent2 color (orange).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shape = $s.Shapes.Item("Elipse 3")

$shape.Fill.Visible = $true
$shape.Fill.Solid()
$shape.Fill.ForeColor.ObjectThemeColor = 6  # msoThemeColorAccent2 -> <a:schemeClr val="accent2"/>
